# Workbook "P_aeropuertos_pasaje" edit:
#  1) Rename the "Mes"/airport-name table headers (row 5) so the leading /
#     double spaces are stripped and "PuertoVallarta" / "SantaLucia" gain
#     their missing space.
#  2) Replace the numeric month value in column C (rows 6-85) with the
#     Spanish month abbreviation ("Ago.", "Jul.", ... ) as text, matching
#     the values already used for the "Mes" column of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames (row 5 of the "Tabla3" table) -----------------------
$ws.Range("G5").Value2 = "Guadalajara"
$ws.Range("H5").Value2 = "Monterrey"
$ws.Range("I5").Value2 = "Tijuana"
$ws.Range("K5").Value2 = "Puerto Vallarta"
$ws.Range("L5").Value2 = "Mérida"
$ws.Range("M5").Value2 = "Guanajuato"
$ws.Range("N5").Value2 = "Culiacán"
$ws.Range("O5").Value2 = "Santa Lucia"
$ws.Range("P5").Value2 = "Otros"

# --- 2) Column C ("Mes"): numeric month -> Spanish month abbreviation ------
$months = @("Ene.", "Feb.", "Mar.", "Abr.", "May.", "Jun.", "Jul.", "Ago.", "Sep.", "Oct.", "Nov.", "Dic.")

for ($row = 6; $row -le 85; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $monthNumber = $cell.Value2
    if ($monthNumber -is [string]) {
        continue
    }
    $idx = [int]$monthNumber - 1
    $cell.Value2 = $months[$idx]
}
